# Katalon UM-Data.xlsx bootstrap update
#
# Commit message: "Changed UM Create User Password to temp Password Policy
# of Upper case, lower case and Number."
#
# The old temp passwords "hello5555" / "hello6666" (all lower-case) are
# replaced everywhere by "Hello5555" (upper + lower + number) on the
# "CreateUser" and "CreateUserErrors" sheets (Password / ConfirmPassword
# columns).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. CreateUser sheet - Password (I) / ConfirmPassword (J) columns
# ---------------------------------------------------------------------------
$wsCreateUser = $wb.Worksheets.Item("CreateUser")
$wsCreateUser.Range("I2").Value = "Hello5555"
$wsCreateUser.Range("J2").Value = "Hello5555"
$wsCreateUser.Range("I3").Value = "Hello5555"
$wsCreateUser.Range("J3").Value = "Hello5555"
$wsCreateUser.Range("I4").Value = "Hello5555"
$wsCreateUser.Range("J4").Value = "Hello5555"

# ---------------------------------------------------------------------------
# 2. CreateUserErrors sheet - Password (I) / ConfirmPassword (J) columns
# ---------------------------------------------------------------------------
$wsCreateUserErrors = $wb.Worksheets.Item("CreateUserErrors")
$wsCreateUserErrors.Range("I2").Value = "Hello5555"
$wsCreateUserErrors.Range("J2").Value = "Hello5555"
$wsCreateUserErrors.Range("I3").Value = "Hello5555"
$wsCreateUserErrors.Range("J3").Value = "Hello5555"
$wsCreateUserErrors.Range("I4").Value = "Hello5555"
$wsCreateUserErrors.Range("J4").Value = "Hello5555"
$wsCreateUserErrors.Range("I5").Value = "Hello5555"
$wsCreateUserErrors.Range("J5").Value = "Hello5555"

# ---------------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping to mirror the saved workbook view
# ---------------------------------------------------------------------------
$wsTestForTyler = $wb.Worksheets.Item("TestForTyler")
$wsTestForTyler.Range("G25").Select() | Out-Null

$wsCreateUser.Range("J8").Select() | Out-Null

# CreateUserErrors ends up being the active / tab-selected sheet, with J5
# selected.
$wsCreateUserErrors.Activate() | Out-Null
$wsCreateUserErrors.Range("J5").Select() | Out-Null
